$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (old rows 2-7 shift down to 3-8).
$ws.Rows.Item(2).Insert()

# The newly inserted row 2 is populated with the same pattern as the row
# that is now row 3 (the old row 2), so copy that row's values down into
# the new row 2.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(2, $col).Value2 = $ws.Cells.Item(3, $col).Value2
}

# ...except column A, which starts the sequence over at 0 instead of 5.
$ws.Cells.Item(2, 1).Value2 = 0

# Update the selection to match the saved workbook state.
$ws.Range("C11").Select()
